$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 688.125
$ws.Range("I28").Value = 688.125
$ws.Range("K28").Value = 688.125
$ws.Range("M28").Value = -203.125
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 5625.1816
$ws.Range("I74").Value = 1968.25
$ws.Range("J74").Value = 7714.857
$ws.Range("K74").Value = 1968.25
$ws.Range("L74").Value = 7714.857
$ws.Range("M74").Value = -1032.25
$ws.Range("N74").Value = -9586.857
$ws.Range("H77").Value = 5625.1816
$ws.Range("I77").Value = 1968.25
$ws.Range("J77").Value = 7714.857
$ws.Range("K77").Value = 9841.25
$ws.Range("L77").Value = 38574.285
$ws.Range("M77").Value = -5161.25
$ws.Range("N77").Value = -47934.285
$ws.Range("H80").Value = 2601.0527
$ws.Range("I80").Value = 2198.3
$ws.Range("J80").Value = 3048.5557
$ws.Range("K80").Value = 6594.900000000001
$ws.Range("L80").Value = 9145.667099999999
$ws.Range("M80").Value = -5596.900000000001
$ws.Range("N80").Value = -11141.6671
$ws.Range("H83").Value = 2601.0527
$ws.Range("I83").Value = 2198.3
$ws.Range("J83").Value = 3048.5557
$ws.Range("K83").Value = 19784.7
$ws.Range("L83").Value = 27437.0013
$ws.Range("M83").Value = -14792.7
$ws.Range("N83").Value = -37421.0013
$ws.Range("H98").Value = 1339.5769
$ws.Range("J98").Value = 4521.5
$ws.Range("L98").Value = 4521.5
$ws.Range("N98").Value = -7517.5
$ws.Range("H112").Value = 1885.2667
$ws.Range("I112").Value = 1145
$ws.Range("J112").Value = 1999.1538
$ws.Range("K112").Value = 3435
$ws.Range("L112").Value = 5997.4614
$ws.Range("M112").Value = -2327
$ws.Range("N112").Value = -8213.4614
$ws.Range("H113").Value = 45459.875
$ws.Range("I113").Value = 102574.9
$ws.Range("J113").Value = 4663.4287
$ws.Range("K113").Value = 102574.9
$ws.Range("L113").Value = 4663.4287
$ws.Range("M113").Value = -99320.89999999999
$ws.Range("N113").Value = -11171.4287
$ws.Range("H122").Value = 1339.5769
$ws.Range("J122").Value = 4521.5
$ws.Range("L122").Value = 13564.5
$ws.Range("N122").Value = -18464.5
$ws.Range("H131").Value = 1200
$ws.Range("I131").Value = 1200
$ws.Range("K131").Value = 3600
$ws.Range("M131").Value = 1440
$ws.Range("H137").Value = 2867.3928
$ws.Range("I137").Value = 1690.5834
$ws.Range("J137").Value = 3750
$ws.Range("K137").Value = 5071.7502
$ws.Range("L137").Value = 11250
$ws.Range("M137").Value = -2521.7502
$ws.Range("N137").Value = -16350

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4307.616
$ws.Range("I32").Value = 3193.0386
$ws.Range("K32").Value = 3193.0386
$ws.Range("M32").Value = -2906.0386
$ws.Range("H63").Value = 4344.25
$ws.Range("I63").Value = 2459
$ws.Range("K63").Value = 2459
$ws.Range("M63").Value = -1773
$ws.Range("H66").Value = 4344.25
$ws.Range("I66").Value = 2459
$ws.Range("K66").Value = 12295
$ws.Range("M66").Value = -8863
$ws.Range("H74").Value = 2533.6155
$ws.Range("I74").Value = 1200
$ws.Range("J74").Value = 3126.3333
$ws.Range("K74").Value = 1200
$ws.Range("L74").Value = 3126.3333
$ws.Range("M74").Value = -326
$ws.Range("N74").Value = -4874.3333
$ws.Range("H77").Value = 2533.6155
$ws.Range("I77").Value = 1200
$ws.Range("J77").Value = 3126.3333
$ws.Range("K77").Value = 6000
$ws.Range("L77").Value = 15631.6665
$ws.Range("M77").Value = -1632
$ws.Range("N77").Value = -24367.6665

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15631113
$ws.Range("I20").Value = 17247582
$ws.Range("J20").Value = 5238.6665
$ws.Range("K20").Value = 17247582
$ws.Range("L20").Value = 5238.6665
$ws.Range("M20").Value = -17247335
$ws.Range("N20").Value = -5732.6665
$ws.Range("H22").Value = 1213
$ws.Range("I22").Value = 1355.4286
$ws.Range("K22").Value = 1355.4286
$ws.Range("M22").Value = -1182.4286
$ws.Range("H68").Value = 23147.5
$ws.Range("J68").Value = 23147.5
$ws.Range("L68").Value = 23147.5
$ws.Range("N68").Value = -24769.5
$ws.Range("H71").Value = 23147.5
$ws.Range("J71").Value = 23147.5
$ws.Range("L71").Value = 69442.5
$ws.Range("N71").Value = -77554.5
$ws.Range("H81").Value = 59233.332
$ws.Range("J81").Value = 59233.332
$ws.Range("L81").Value = 59233.332
$ws.Range("N81").Value = -61355.332
$ws.Range("H84").Value = 59233.332
$ws.Range("J84").Value = 59233.332
$ws.Range("L84").Value = 177699.996
$ws.Range("N84").Value = -188307.996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24317.826
$ws.Range("I31").Value = 29166.889
$ws.Range("J31").Value = 6861.2
$ws.Range("K31").Value = 29166.889
$ws.Range("L31").Value = 6861.2
$ws.Range("M31").Value = -28871.889
$ws.Range("N31").Value = -7451.2
$ws.Range("H34").Value = 24317.826
$ws.Range("I34").Value = 29166.889
$ws.Range("J34").Value = 6861.2
$ws.Range("K34").Value = 29166.889
$ws.Range("L34").Value = 6861.2
$ws.Range("M34").Value = -28964.889
$ws.Range("N34").Value = -7265.2
$ws.Range("H68").Value = 20663.334
$ws.Range("J68").Value = 20663.334
$ws.Range("L68").Value = 20663.334
$ws.Range("N68").Value = -22161.334
$ws.Range("H71").Value = 20663.334
$ws.Range("J71").Value = 20663.334
$ws.Range("L71").Value = 61990.00199999999
$ws.Range("N71").Value = -69478.00199999999
$ws.Range("H132").Value = 3124.7144
$ws.Range("I132").Value = 3158.0715
$ws.Range("K132").Value = 9474.2145
$ws.Range("M132").Value = -6944.2145
$ws.Range("H141").Value = 200967.75
$ws.Range("J141").Value = 200967.75
$ws.Range("L141").Value = 200967.75
$ws.Range("N141").Value = -211327.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 6340.5
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H122").Value = 836
$ws.Range("I122").Value = 840
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 7560
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -5110
$ws.Range("N122").Value = -12100
$ws.Range("H131").Value = 1987.6538
$ws.Range("I131").Value = 1993.3334
$ws.Range("K131").Value = 5980.0002
$ws.Range("M131").Value = -940.0002000000004
$ws.Range("H140").Value = 2776.6
$ws.Range("I140").Value = 2630.6924
$ws.Range("K140").Value = 7892.0772
$ws.Range("M140").Value = -2712.0772

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 847.55884
$ws.Range("J2").Value = 1303.6
$ws.Range("L2").Value = 1303.6
$ws.Range("N2").Value = -1529.6
$ws.Range("H102").Value = 2153.842
$ws.Range("I102").Value = 2153.842
$ws.Range("K102").Value = 2153.842
$ws.Range("M102").Value = -531.8420000000001
$ws.Range("H132").Value = 3822.6
$ws.Range("I132").Value = 3025
$ws.Range("K132").Value = 9075
$ws.Range("M132").Value = -6545

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 32687.104
$ws.Range("I7").Value = 33679.668
$ws.Range("K7").Value = 33679.668
$ws.Range("M7").Value = -33567.668
$ws.Range("H16").Value = 728.5833
$ws.Range("I16").Value = 749.36365
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 749.36365
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -579.36365
$ws.Range("N16").Value = -840
$ws.Range("H46").Value = 2981.6333
$ws.Range("I46").Value = 2772.75
$ws.Range("J46").Value = 3399.4
$ws.Range("K46").Value = 2772.75
$ws.Range("L46").Value = 3399.4
$ws.Range("M46").Value = -2584.75
$ws.Range("N46").Value = -3775.4
$ws.Range("H100").Value = 4367.8213
$ws.Range("I100").Value = 2712.0588
$ws.Range("K100").Value = 2712.0588
$ws.Range("M100").Value = -2171.0588
$ws.Range("H122").Value = 95135.37
$ws.Range("I122").Value = 128776.375
$ws.Range("J122").Value = 5426
$ws.Range("K122").Value = 386329.125
$ws.Range("L122").Value = 16278
$ws.Range("M122").Value = -383879.125
$ws.Range("N122").Value = -21178
$ws.Range("H123").Value = 77214.5
$ws.Range("J123").Value = 77214.5
$ws.Range("L123").Value = 77214.5
$ws.Range("N123").Value = -87014.5
$ws.Range("H126").Value = 32687.104
$ws.Range("I126").Value = 33679.668
$ws.Range("K126").Value = 101039.004
$ws.Range("M126").Value = -98569.00399999999
$ws.Range("H132").Value = 4176.8613
$ws.Range("I132").Value = 3275.6553
$ws.Range("K132").Value = 9826.965899999999
$ws.Range("M132").Value = -7296.965899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H62").Value = 19285
$ws.Range("I62").Value = 4998.75
$ws.Range("J62").Value = 38333.332
$ws.Range("K62").Value = 4998.75
$ws.Range("L62").Value = 38333.332
$ws.Range("M62").Value = -4374.75
$ws.Range("N62").Value = -39581.332
$ws.Range("H65").Value = 19285
$ws.Range("I65").Value = 4998.75
$ws.Range("J65").Value = 38333.332
$ws.Range("K65").Value = 24993.75
$ws.Range("L65").Value = 191666.66
$ws.Range("M65").Value = -21873.75
$ws.Range("N65").Value = -197906.66
$ws.Range("H126").Value = 21202.084
$ws.Range("I126").Value = 34231.145
$ws.Range("J126").Value = 2961.4
$ws.Range("K126").Value = 102693.435
$ws.Range("L126").Value = 8884.200000000001
$ws.Range("M126").Value = -100223.435
$ws.Range("N126").Value = -13824.2

Write-Host "Applied Leviathan_Profits updates"